$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text (e.g. "1.00", "577.65").
# Force Text format on the D cells we are about to rewrite so COM/Excel does
# not silently coerce them into numbers (which would drop the literal
# formatting, e.g. "1.00" -> 1). Each contiguous run is a separate statement
# since this host does not honour comma multi-area Range refs.
$ws.Range("D2:D3").NumberFormat = "@"
$ws.Range("D5:D14").NumberFormat = "@"
$ws.Range("D16:D22").NumberFormat = "@"
$ws.Range("D25:D26").NumberFormat = "@"
$ws.Range("D28:D32").NumberFormat = "@"
$ws.Range("D34:D38").NumberFormat = "@"
$ws.Range("D40:D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

$ws.Range("D2").Value = '64.912.15'
$ws.Range("E2").Value = '  -1.21%  '

$ws.Range("D3").Value = '3.389.84'
$ws.Range("E3").Value = '  -3.99%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '577.65'
$ws.Range("E5").Value = '  -4.05%  '

$ws.Range("D6").Value = '134.98'
$ws.Range("E6").Value = '  -6.09%  '

$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = '3.388.06'
$ws.Range("E8").Value = '  -4.03%  '

$ws.Range("D9").Value = '0.491'
$ws.Range("E9").Value = '  -2.15%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.120'
$ws.Range("E10").Value = '  -10.28%  '

$ws.Range("B11").Value = 'Toncoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D11").Value = '7.06'
$ws.Range("E11").Value = '  -10.02%  '

$ws.Range("D12").Value = '0.370'
$ws.Range("E12").Value = '  -8.19%  '

$ws.Range("D13").Value = '3.966.04'
$ws.Range("E13").Value = '  -4.04%  '

$ws.Range("D14").Value = '0.0000176'
$ws.Range("E14").Value = '  -11.28%  '

$ws.Range("E15").Value = '  -1.84%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.396.04'
$ws.Range("E16").Value = '  -3.63%  '

$ws.Range("D17").Value = '64.900.63'
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '25.95'
$ws.Range("E18").Value = '  -8.75%  '

$ws.Range("D19").Value = '9.45'
$ws.Range("E19").Value = '  -14.32%  '

$ws.Range("D20").Value = '5.79'
$ws.Range("E20").Value = '  -6.63%  '

$ws.Range("D21").Value = '13.44'
$ws.Range("E21").Value = '  -6.01%  '

$ws.Range("D22").Value = '378.66'
$ws.Range("E22").Value = '  -8.83%  '

$ws.Range("E23").Value = '  -8.39%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("D25").Value = '71.74'
$ws.Range("E25").Value = '  -7.47%  '

$ws.Range("D26").Value = '3.525.58'
$ws.Range("E26").Value = '  -4.02%  '

$ws.Range("E27").Value = '  -10.62%  '

$ws.Range("D28").Value = '1.00'
$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").Value = '2.18'
$ws.Range("E29").Value = '  -10.49%  '

$ws.Range("D30").Value = '6.95'
$ws.Range("E30").Value = '  -10.01%  '

$ws.Range("D31").Value = '7.93'
$ws.Range("E31").Value = '  -10.53%  '

$ws.Range("D32").Value = '3.399.73'
$ws.Range("E32").Value = '  -3.74%  '

$ws.Range("E33").Value = '  +0.01%  '

$ws.Range("D34").Value = '0.141'
$ws.Range("E34").Value = '  -7.57%  '

$ws.Range("D35").Value = '22.66'
$ws.Range("E35").Value = '  -6.94%  '

$ws.Range("D36").Value = '168.39'
$ws.Range("E36").Value = '  -3.32%  '

$ws.Range("D37").Value = '6.61'
$ws.Range("E37").Value = '  -12.28%  '

$ws.Range("D38").Value = '1.13'
$ws.Range("E38").Value = '  -11.71%  '

$ws.Range("E39").Value = '  -7.66%  '

$ws.Range("D40").Value = '4.62'
$ws.Range("E40").Value = '  -12.22%  '

$ws.Range("D41").Value = '0.0747'
$ws.Range("E41").Value = '  -8.69%  '

$ws.Range("D42").Value = '0.807'
$ws.Range("E42").Value = '  -6.08%  '

$ws.Range("D43").Value = '43.32'
$ws.Range("E43").Value = '  -4.65%  '

$ws.Range("E44").Value = '  +0.14%  '

$ws.Range("D45").Value = '4.29'
$ws.Range("E45").Value = '  -15.30%  '

$ws.Range("E46").Value = '  -10.42%  '

$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("E48").Value = '  -6.21%  '

$ws.Range("D49").Value = '6.41'
$ws.Range("E49").Value = '  -8.35%  '

$ws.Range("D50").Value = '2.01'
$ws.Range("E50").Value = '  -15.08%  '

$ws.Range("D51").Value = '2.147.43'
$ws.Range("E51").Value = '  -8.92%  '
